$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, pushing all following rows down by one.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row with the new descriptor entry.
$ws.Cells.Item(32, 1).Value = "Tesla Model X I · Рестайлинг"
$ws.Cells.Item(32, 2).Value = "https://cars.av.by/tesla/model-x/109287918"
$ws.Cells.Item(32, 3).Value = 270985

# The insert pushed a duplicate of the old last row (64) down to row 65;
# remove it so the sheet still ends at row 64.
$ws.Rows.Item(65).Delete()
